$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (the "Förändrad" date column) from row 2 through row 43:
# bump the serial date value from 45843 to 45844 (2025-07-05 -> 2025-07-06).
for ($row = 2; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45843) {
        $cell.Value = 45844
    }
}
